$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 41; existing rows 41-51 shift down to 42-52.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with its data.
$ws.Cells.Item(41, 1).Value = 7
$ws.Cells.Item(41, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(41, 3).Value = "Ñuble"
$ws.Cells.Item(41, 4).Value = 45202
$ws.Cells.Item(41, 5).Value = 16
$ws.Cells.Item(41, 6).Value = 300000000
$ws.Cells.Item(41, 7).Value = "Espárragos"
$ws.Cells.Item(41, 8).Value = "Sin especificar"
$ws.Cells.Item(41, 9).Value = "Primera"
$ws.Cells.Item(41, 10).Value = 300
$ws.Cells.Item(41, 11).Value = 1500
$ws.Cells.Item(41, 12).Value = 1500
$ws.Cells.Item(41, 13).Value = 1500
$ws.Cells.Item(41, 14).Value = "$/kilo"
$ws.Cells.Item(41, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(41, 16).Value = 1500
$ws.Cells.Item(41, 17).Value = 1
$ws.Cells.Item(41, 18).Value = "Hortaliza"

# Row that was previously 50 is now 51 after the insert; only its
# "Origen" column changes (Región de Ñuble -> Provincia de Diguillín).
$ws.Cells.Item(51, 15).Value = "Provincia de Diguillín"
